$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 40766.668
$ws.Range("J21").Value = 52150
$ws.Range("L21").Value = 52150
$ws.Range("N21").Value = -53086

$ws.Range("H23").Value = 40766.668
$ws.Range("J23").Value = 52150
$ws.Range("L23").Value = 52150
$ws.Range("N23").Value = -52618

$ws.Range("H40").Value = 1740
$ws.Range("I40").Value = 1698
$ws.Range("K40").Value = 1698
$ws.Range("M40").Value = -1523

$ws.Range("H92").Value = 1942.05
$ws.Range("I92").Value = 1851
$ws.Range("K92").Value = 1851
$ws.Range("M92").Value = -603

$ws.Range("H106").Value = 2505002.5
$ws.Range("I106").Value = 2505002.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2505002.5
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -2504371.5

$ws.Range("H118").Value = 3344.75
$ws.Range("I118").Value = 179
$ws.Range("J118").Value = 4400
$ws.Range("K118").Value = 537
$ws.Range("L118").Value = 13200
$ws.Range("M118").Value = 1120
$ws.Range("N118").Value = -16514

$ws.Range("H137").Value = 4167557.2
$ws.Range("I137").Value = 854.7805
$ws.Range("K137").Value = 2564.3415
$ws.Range("M137").Value = -14.3415

$ws.Range("H139").Value = 52354
$ws.Range("J139").Value = 52354
$ws.Range("L139").Value = 52354
$ws.Range("N139").Value = -62634

$ws.Range("H140").Value = 49933.332
$ws.Range("J140").Value = 49933.332
$ws.Range("L140").Value = 49933.332
$ws.Range("N140").Value = -60293.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12501351
$ws.Range("I61").Value = 16130459
$ws.Range("K61").Value = 16130459
$ws.Range("M61").Value = -16130247

$ws.Range("H74").Value = 13516603
$ws.Range("I74").Value = 19232382
$ws.Range("K74").Value = 19232382
$ws.Range("M74").Value = -19231508

$ws.Range("H77").Value = 13516603
$ws.Range("I77").Value = 19232382
$ws.Range("K77").Value = 96161910
$ws.Range("M77").Value = -96157542

$ws.Range("H132").Value = 13892676
$ws.Range("I132").Value = 41671492
$ws.Range("J132").Value = 3268.3333
$ws.Range("K132").Value = 125014476
$ws.Range("L132").Value = 9804.999899999999
$ws.Range("M132").Value = -125011946
$ws.Range("N132").Value = -14864.9999

$ws.Range("H136").Value = 12501351
$ws.Range("I136").Value = 16130459
$ws.Range("K136").Value = 48391377
$ws.Range("M136").Value = -48388827

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 5936
$ws.Range("J38").Value = 5936
$ws.Range("L38").Value = 5936
$ws.Range("N38").Value = -6768

$ws.Range("H92").Value = 14000
$ws.Range("J92").Value = 14000
$ws.Range("L92").Value = 14000
$ws.Range("N92").Value = -18992

$ws.Range("H134").Value = 2804.75
$ws.Range("I134").Value = 1812.0605
$ws.Range("J134").Value = 5782.8184
$ws.Range("K134").Value = 5436.181500000001
$ws.Range("L134").Value = 17348.4552
$ws.Range("M134").Value = -2901.181500000001
$ws.Range("N134").Value = -22418.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851152.5
$ws.Range("I31").Value = 3346.426
$ws.Range("J31").Value = 111111656
$ws.Range("K31").Value = 3346.426
$ws.Range("L31").Value = 111111656
$ws.Range("M31").Value = -3051.426
$ws.Range("N31").Value = -111112246

$ws.Range("H34").Value = 5851152.5
$ws.Range("I34").Value = 3346.426
$ws.Range("J34").Value = 111111656
$ws.Range("K34").Value = 3346.426
$ws.Range("L34").Value = 111111656
$ws.Range("M34").Value = -3144.426
$ws.Range("N34").Value = -111112060

$ws.Range("H58").Value = 1458.449
$ws.Range("I58").Value = 758.62964
$ws.Range("J58").Value = 2317.318
$ws.Range("K58").Value = 758.62964
$ws.Range("L58").Value = 2317.318
$ws.Range("M58").Value = -555.62964
$ws.Range("N58").Value = -2723.318

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = 0

$ws.Range("H132").Value = 10205605
$ws.Range("I132").Value = 13514734
$ws.Range("J132").Value = 2456.5
$ws.Range("K132").Value = 40544202
$ws.Range("L132").Value = 7369.5
$ws.Range("M132").Value = -40541672
$ws.Range("N132").Value = -12429.5

$ws.Range("H134").Value = 1311.6666
$ws.Range("I134").Value = 1360.4333
$ws.Range("J134").Value = 824
$ws.Range("K134").Value = 4081.2999
$ws.Range("L134").Value = 2472
$ws.Range("M134").Value = -1546.2999
$ws.Range("N134").Value = -7542

$ws.Range("H136").Value = 1458.449
$ws.Range("I136").Value = 758.62964
$ws.Range("J136").Value = 2317.318
$ws.Range("K136").Value = 2275.88892
$ws.Range("L136").Value = 6951.954000000001
$ws.Range("M136").Value = 274.1110800000001
$ws.Range("N136").Value = -12051.954

$ws.Range("H140").Value = 38769.75
$ws.Range("J140").Value = 38769.75
$ws.Range("L140").Value = 38769.75
$ws.Range("N140").Value = -49129.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5776.96
$ws.Range("I132").Value = 4865.3335
$ws.Range("K132").Value = 14596.0005
$ws.Range("M132").Value = -12066.0005

$ws.Range("H134").Value = 35795.6
$ws.Range("J134").Value = 35795.6
$ws.Range("L134").Value = 107386.8
$ws.Range("N134").Value = -112456.8

$ws.Range("H138").Value = 59249.25
$ws.Range("J138").Value = 59249.25
$ws.Range("L138").Value = 59249.25
$ws.Range("N138").Value = -69529.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 686.8333
$ws.Range("J46").Value = 888.3
$ws.Range("L46").Value = 888.3
$ws.Range("N46").Value = -1264.3

$ws.Range("H74").Value = 25000
$ws.Range("J74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("N74").Value = -26996

$ws.Range("H77").Value = 25000
$ws.Range("J77").Value = 25000
$ws.Range("L77").Value = 75000
$ws.Range("N77").Value = -84984

$ws.Range("H101").Value = 8439.143
$ws.Range("J101").Value = 8439.143
$ws.Range("L101").Value = 8439.143
$ws.Range("N101").Value = -14929.143

$ws.Range("H132").Value = 7698119
$ws.Range("I132").Value = 3416.465
$ws.Range("J132").Value = 22737764
$ws.Range("K132").Value = 10249.395
$ws.Range("L132").Value = 68213292
$ws.Range("M132").Value = -7719.395
$ws.Range("N132").Value = -68218352

$ws.Range("H136").Value = 10642546
$ws.Range("I136").Value = 15153176
$ws.Range("J136").Value = 10347.5
$ws.Range("K136").Value = 45459528
$ws.Range("L136").Value = 31042.5
$ws.Range("M136").Value = -45456978
$ws.Range("N136").Value = -36142.5

$ws.Range("H139").Value = 58893
$ws.Range("J139").Value = 58893
$ws.Range("L139").Value = 58893
$ws.Range("N139").Value = -69173

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 39800
$ws.Range("J42").Value = 39800
$ws.Range("L42").Value = 39800
$ws.Range("N42").Value = -40556

$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H80").Value = 37843.285
$ws.Range("J80").Value = 37483.832
$ws.Range("L80").Value = 37483.832
$ws.Range("N80").Value = -39479.832

$ws.Range("H81").Value = 941
$ws.Range("I81").Value = 980
$ws.Range("J81").Value = 902
$ws.Range("K81").Value = 1960
$ws.Range("L81").Value = 1804
$ws.Range("M81").Value = -899
$ws.Range("N81").Value = -3926

$ws.Range("H83").Value = 37843.285
$ws.Range("J83").Value = 37483.832
$ws.Range("L83").Value = 112451.496
$ws.Range("N83").Value = -122435.496

$ws.Range("H84").Value = 941
$ws.Range("I84").Value = 980
$ws.Range("J84").Value = 902
$ws.Range("K84").Value = 9800
$ws.Range("L84").Value = 9020
$ws.Range("M84").Value = -4496
$ws.Range("N84").Value = -19628

$ws.Range("H132").Value = 3084.3333
$ws.Range("I132").Value = 2014
$ws.Range("J132").Value = 4307.5713
$ws.Range("K132").Value = 6042
$ws.Range("L132").Value = 12922.7139
$ws.Range("M132").Value = -3512
$ws.Range("N132").Value = -17982.7139

$ws.Range("H136").Value = 1168.4117
$ws.Range("I136").Value = 1012.37036
$ws.Range("J136").Value = 1770.2858
$ws.Range("K136").Value = 3037.11108
$ws.Range("L136").Value = 5310.857400000001
$ws.Range("M136").Value = -487.1110800000001
$ws.Range("N136").Value = -10410.8574
